$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 282.08334
$ws.Range("I33").Value = 153.88889
$ws.Range("K33").Value = 153.88889
$ws.Range("M33").Value = 75.11111

$ws.Range("H76").Value = 13300
$ws.Range("I76").Value = 17285.715
$ws.Range("K76").Value = 17285.715
$ws.Range("M76").Value = -16970.715

$ws.Range("H79").Value = 13300
$ws.Range("I79").Value = 17285.715
$ws.Range("K79").Value = 17285.715
$ws.Range("M79").Value = -16193.715

$ws.Range("H86").Value = 31258576
$ws.Range("I86").Value = 13333.889
$ws.Range("J86").Value = 71431030
$ws.Range("K86").Value = 13333.889
$ws.Range("L86").Value = 71431030
$ws.Range("M86").Value = -12210.889
$ws.Range("N86").Value = -71433276

$ws.Range("H89").Value = 31258576
$ws.Range("I89").Value = 13333.889
$ws.Range("J89").Value = 71431030
$ws.Range("K89").Value = 66669.44499999999
$ws.Range("L89").Value = 357155150
$ws.Range("M89").Value = -61053.44499999999
$ws.Range("N89").Value = -357166382

$ws.Range("H113").Value = 29415736
$ws.Range("I113").Value = 83336070
$ws.Range("J113").Value = 4644.364
$ws.Range("K113").Value = 83336070
$ws.Range("L113").Value = 4644.364
$ws.Range("M113").Value = -83332816
$ws.Range("N113").Value = -11152.364

$ws.Range("H128").Value = 14830.368
$ws.Range("J128").Value = 14830.368
$ws.Range("L128").Value = 14830.368
$ws.Range("N128").Value = -24790.368

$ws.Range("H130").Value = 20605.715
$ws.Range("J130").Value = 20605.715
$ws.Range("L130").Value = 20605.715
$ws.Range("N130").Value = -30645.715

$ws.Range("H132").Value = 2026.8823
$ws.Range("I132").Value = 1351.675
$ws.Range("J132").Value = 4482.1816
$ws.Range("K132").Value = 4055.025
$ws.Range("L132").Value = 13446.5448
$ws.Range("M132").Value = -1525.025
$ws.Range("N132").Value = -18506.5448

$ws.Range("H135").Value = 21890.041
$ws.Range("I135").Value = 28642.805
$ws.Range("K135").Value = 257785.245
$ws.Range("M135").Value = -255250.245

$ws.Range("H141").Value = 1605.5
$ws.Range("I141").Value = 1177.3823
$ws.Range("J141").Value = 3425
$ws.Range("K141").Value = 3532.1469
$ws.Range("L141").Value = 10275
$ws.Range("M141").Value = 1647.8531
$ws.Range("N141").Value = -20635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23130.014
$ws.Range("I32").Value = 5303.8623
$ws.Range("K32").Value = 5303.8623
$ws.Range("M32").Value = -5016.8623

$ws.Range("H45").Value = 2534.9167
$ws.Range("I45").Value = 2128.625
$ws.Range("J45").Value = 2859.95
$ws.Range("K45").Value = 2128.625
$ws.Range("L45").Value = 2859.95
$ws.Range("M45").Value = -1751.625
$ws.Range("N45").Value = -3613.95

$ws.Range("H122").Value = 2828.5715
$ws.Range("I122").Value = 2828.5715
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8485.7145
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6035.7145
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15878323
$ws.Range("I20").Value = 22228254
$ws.Range("J20").Value = 3498.3333
$ws.Range("K20").Value = 22228254
$ws.Range("L20").Value = 3498.3333
$ws.Range("M20").Value = -22228007
$ws.Range("N20").Value = -3992.3333

$ws.Range("H134").Value = 1976.3914
$ws.Range("I134").Value = 1865.8536
$ws.Range("J134").Value = 2882.8
$ws.Range("K134").Value = 5597.560799999999
$ws.Range("L134").Value = 8648.400000000001
$ws.Range("M134").Value = -3062.560799999999
$ws.Range("N134").Value = -13718.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1795.8788
$ws.Range("I31").Value = 1498.8
$ws.Range("J31").Value = 2252.923
$ws.Range("K31").Value = 1498.8
$ws.Range("L31").Value = 2252.923
$ws.Range("M31").Value = -1203.8
$ws.Range("N31").Value = -2842.923

$ws.Range("H34").Value = 1795.8788
$ws.Range("I34").Value = 1498.8
$ws.Range("J34").Value = 2252.923
$ws.Range("K34").Value = 1498.8
$ws.Range("L34").Value = 2252.923
$ws.Range("M34").Value = -1296.8
$ws.Range("N34").Value = -2656.923

$ws.Range("H58").Value = 792.87036
$ws.Range("I58").Value = 780.1667
$ws.Range("J58").Value = 894.5
$ws.Range("K58").Value = 780.1667
$ws.Range("L58").Value = 894.5
$ws.Range("M58").Value = -577.1667
$ws.Range("N58").Value = -1300.5

$ws.Range("H122").Value = 1995.1333
$ws.Range("J122").Value = 3169
$ws.Range("L122").Value = 9507
$ws.Range("N122").Value = -14407

$ws.Range("H132").Value = 2935.5625
$ws.Range("I132").Value = 2856
$ws.Range("J132").Value = 3110.6
$ws.Range("K132").Value = 8568
$ws.Range("L132").Value = 9331.799999999999
$ws.Range("M132").Value = -6038
$ws.Range("N132").Value = -14391.8

$ws.Range("H136").Value = 792.87036
$ws.Range("I136").Value = 780.1667
$ws.Range("J136").Value = 894.5
$ws.Range("K136").Value = 2340.5001
$ws.Range("L136").Value = 2683.5
$ws.Range("M136").Value = 209.4998999999998
$ws.Range("N136").Value = -7783.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1334266.8
$ws.Range("I5").Value = 847.7143
$ws.Range("J5").Value = 4445577.5
$ws.Range("K5").Value = 2543.1429
$ws.Range("L5").Value = 13336732.5
$ws.Range("M5").Value = -2431.1429
$ws.Range("N5").Value = -13336956.5

$ws.Range("H129").Value = 1918.238
$ws.Range("J129").Value = 2408.7856
$ws.Range("L129").Value = 7226.3568
$ws.Range("N129").Value = -17226.3568

$ws.Range("H135").Value = 1334266.8
$ws.Range("I135").Value = 847.7143
$ws.Range("J135").Value = 4445577.5
$ws.Range("K135").Value = 7629.428699999999
$ws.Range("L135").Value = 40010197.5
$ws.Range("M135").Value = -5094.428699999999
$ws.Range("N135").Value = -40015267.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 9890
$ws.Range("J93").Value = 9890
$ws.Range("L93").Value = 9890
$ws.Range("N93").Value = -13634

$ws.Range("H102").Value = 1704.1666
$ws.Range("I102").Value = 1415
$ws.Range("K102").Value = 1415
$ws.Range("M102").Value = 207

$ws.Range("H122").Value = 3901
$ws.Range("I122").Value = 5002.3335
$ws.Range("J122").Value = 3075
$ws.Range("K122").Value = 15007.0005
$ws.Range("L122").Value = 9225
$ws.Range("M122").Value = -12557.0005
$ws.Range("N122").Value = -14125

$ws.Range("H130").Value = 47086
$ws.Range("J130").Value = 47086
$ws.Range("L130").Value = 47086
$ws.Range("N130").Value = -57126

$ws.Range("H135").Value = 53641.117
$ws.Range("J135").Value = 53641.117
$ws.Range("L135").Value = 53641.117
$ws.Range("N135").Value = -63781.117

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2133.3333
$ws.Range("I40").Value = 2200
$ws.Range("K40").Value = 2200
$ws.Range("M40").Value = -2064

$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622

$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112

$ws.Range("H93").Value = 2278.2222
$ws.Range("I93").Value = 1750
$ws.Range("K93").Value = 1750
$ws.Range("M93").Value = -502

$ws.Range("H122").Value = 3168.889
$ws.Range("I122").Value = 2357.1428
$ws.Range("K122").Value = 7071.428400000001
$ws.Range("M122").Value = -4621.428400000001

$ws.Range("H132").Value = 2503.8845
$ws.Range("I132").Value = 2123.842
$ws.Range("J132").Value = 3535.4285
$ws.Range("K132").Value = 6371.526
$ws.Range("L132").Value = 10606.2855
$ws.Range("M132").Value = -3841.526
$ws.Range("N132").Value = -15666.2855

$ws.Range("H136").Value = 1719.4
$ws.Range("I136").Value = 1805.4
$ws.Range("J136").Value = 1375.4
$ws.Range("K136").Value = 5416.200000000001
$ws.Range("L136").Value = 4126.200000000001
$ws.Range("M136").Value = -2866.200000000001
$ws.Range("N136").Value = -9226.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3555.7778
$ws.Range("I62").Value = 3557.4285
$ws.Range("J62").Value = 3550
$ws.Range("K62").Value = 3557.4285
$ws.Range("L62").Value = 3550
$ws.Range("M62").Value = -2933.4285
$ws.Range("N62").Value = -4798

$ws.Range("H65").Value = 3555.7778
$ws.Range("I65").Value = 3557.4285
$ws.Range("J65").Value = 3550
$ws.Range("K65").Value = 17787.1425
$ws.Range("L65").Value = 17750
$ws.Range("M65").Value = -14667.1425
$ws.Range("N65").Value = -23990

$ws.Range("H93").Value = 25446.334
$ws.Range("J93").Value = 25446.334
$ws.Range("L93").Value = 25446.334
$ws.Range("N93").Value = -30438.334

$ws.Range("H107").Value = 1482.9
$ws.Range("I107").Value = 1049.6364
$ws.Range("J107").Value = 2674.375
$ws.Range("K107").Value = 3148.9092
$ws.Range("L107").Value = 8023.125
$ws.Range("M107").Value = -1228.9092
$ws.Range("N107").Value = -11863.125

$ws.Range("H122").Value = 4807.1
$ws.Range("I122").Value = 5284.684
$ws.Range("J122").Value = 3982.182
$ws.Range("K122").Value = 15854.052
$ws.Range("L122").Value = 11946.546
$ws.Range("M122").Value = -13404.052
$ws.Range("N122").Value = -16846.546

$ws.Range("H132").Value = 4210.409
$ws.Range("I132").Value = 4285.8423
$ws.Range("J132").Value = 3732.6667
$ws.Range("K132").Value = 12857.5269
$ws.Range("L132").Value = 11198.0001
$ws.Range("M132").Value = -10327.5269
$ws.Range("N132").Value = -16258.0001

$ws.Range("H136").Value = 1853.2903
$ws.Range("I136").Value = 1992.3846
$ws.Range("J136").Value = 1130
$ws.Range("K136").Value = 5977.1538
$ws.Range("L136").Value = 3390
$ws.Range("M136").Value = -3427.1538
$ws.Range("N136").Value = -8490
